$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, pushing the existing row 30 (and below) down to row 31.
$ws.Rows.Item(30).Insert()

# Fill the new row 30 with the updated weekly record.
$ws.Cells.Item(30, 1).Value = 2
$ws.Cells.Item(30, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(30, 3).Value = "Coquimbo"
$ws.Cells.Item(30, 4).Value = 44448
$ws.Cells.Item(30, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(30, 5).Value = 4
$ws.Cells.Item(30, 6).Value = 100112022
$ws.Cells.Item(30, 7).Value = "Arveja Verde"
$ws.Cells.Item(30, 8).Value = "Perfection"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 400
$ws.Cells.Item(30, 11).Value = 28000
$ws.Cells.Item(30, 12).Value = 30000
$ws.Cells.Item(30, 13).Value = 29000
$ws.Cells.Item(30, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(30, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(30, 16).Value = 1160
$ws.Cells.Item(30, 17).Value = 25
$ws.Cells.Item(30, 18).Value = "Hortaliza"
